$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = "30.517.83"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).Value = "  +1.66%  "

# Row 3
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = "2.015.63"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Value = "  +5.61%  "

# Row 4
$ws.Cells.Item(4,5).Value = "  +0.00%  "

# Row 5
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "325.12"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = "  +1.61%  "

# Row 6
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "1.002"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = "  +0.02%  "

# Row 7
$ws.Cells.Item(7,5).Value = "  +1.49%  "

# Row 8
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "0.4164"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value = "  +2.88%  "

# Row 9
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "0.08778"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value = "  +5.93%  "

# Row 10
$ws.Cells.Item(10,2).Value = "OKB"
$ws.Cells.Item(10,3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "43.47"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = "  +3.50%  "

# Row 11
$ws.Cells.Item(11,2).Value = "Polygon"
$ws.Cells.Item(11,3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "1.135"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = "  +3.21%  "

# Row 12
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "24.66"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value = "  +2.64%  "

# Row 13
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = "2.018.51"
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value = "  +5.00%  "

# Row 14
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = "6.609"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = "  +3.21%  "

# Row 15
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = "7.492"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value = "  +3.21%  "

# Row 16
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = "1.000"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Value = "  -0.11%  "

# Row 17
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = "94.57"
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Value = "  +2.95%  "

# Row 18
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = "0.00001116"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value = "  +1.86%  "

# Row 19
$ws.Cells.Item(19,5).Value = "  +0.38%  "

# Row 20
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "18.95"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value = "  +4.69%  "

# Row 21
$ws.Cells.Item(21,5).Value = "  +0.02%  "

# Row 22
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "6.225"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = "  +4.90%  "

# Row 23
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = "30.579.33"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value = "  +1.76%  "

# Row 24
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "11.89"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value = "  +5.40%  "

# Row 25
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "2.234"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = "  +2.04%  "

# Row 26
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = "2.246.99"
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value = "  +4.90%  "

# Row 27
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "22.44"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = "  +0.10%  "

# Row 28
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = "163.19"
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).Value = "  +0.58%  "

# Row 29
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "2.433"
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Value = "  +6.53%  "

# Row 30
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = "131.54"
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Value = "  +2.19%  "

# Row 31
$ws.Cells.Item(31,5).Value = "  +2.08%  "

# Row 32
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = "0.1053"
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).Value = "  +1.60%  "

# Row 33
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = "6.102"
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).Value = "  +2.32%  "

# Row 34
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = "3.834"
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value = "  +0.49%  "

# Row 35
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = "1.364"
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).Value = "  +14.39%  "

# Row 36
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = "0.02527"
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).Value = "  +3.53%  "

# Row 37
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = "5.490"
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).Value = "  +2.67%  "

# Row 38
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = "0.06666"
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).Value = "  +4.99%  "

# Row 39
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = "12.32"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value = "  +8.31%  "

# Row 40
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "9.105"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = "  +4.86%  "

# Row 41
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = "0.2200"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value = "  +2.27%  "

# Row 42
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "0.6690"
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).Value = "  +1.20%  "

# Row 43
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "1.233"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value = "  +1.82%  "

# Row 44
$ws.Cells.Item(44,5).Value = "  +0.03%  "

# Row 45
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "13.66"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value = "  +2.33%  "

# Row 46
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = "0.6194"
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value = "  +1.84%  "

# Row 47
$ws.Cells.Item(47,5).Value = "  +0.06%  "

# Row 48
$ws.Cells.Item(48,5).Value = "  +0.88%  "

# Row 49
$ws.Cells.Item(49,5).Value = "  +4.92%  "

# Row 50
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = "124.84"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value = "  +1.46%  "

# Row 51
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = "81.20"
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value = "  +3.41%  "
